# Updates the cryptos list worksheet with refreshed prices / volume(1h)
# percentages, and corrects the ordering of two coin pairs whose ranking
# swapped (Chainlink/BitcoinCash at rows 18-19, Cronos/Algorand at rows 50-51).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => column letter => new value (only cells that actually changed).
$changes = @{
    2  = @{ D = "26.801.32"; E = "  -2.40%  " }
    3  = @{ D = "1.562.04";  E = "  -0.73%  " }
    4  = @{ E = "  +0.12%  " }
    5  = @{ D = "205.61";    E = "  -0.94%  " }
    6  = @{ D = "0.488";     E = "  -1.99%  " }
    7  = @{ E = "  +0.09%  " }
    8  = @{ D = "21.78";     E = "  -2.24%  " }
    9  = @{ D = "0.247";     E = "  -1.14%  " }
    10 = @{ D = "0.0584";    E = "  -1.89%  " }
    11 = @{ D = "0.0866";    E = "  +0.00%  " }
    12 = @{ D = "1.783.09";  E = "  -0.82%  " }
    13 = @{ D = "1.562.09";  E = "  -1.06%  " }
    14 = @{ D = "3.73";      E = "  -2.47%  " }
    15 = @{ D = "0.512";     E = "  -1.41%  " }
    16 = @{ D = "26.832.60"; E = "  -2.31%  " }
    17 = @{ D = "60.92";     E = "  -4.07%  " }
    18 = @{ B = "Chainlink";    C = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link";     D = "7.36";     E = "  +0.62%  " }
    19 = @{ B = "BitcoinCash";  C = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch";     D = "213.92";  E = "  -0.15%  " }
    20 = @{ D = "0.0₃0677";   E = "  -2.25%  " }
    21 = @{ E = "  +0.17%  " }
    22 = @{ D = "4.12";      E = "  -0.69%  " }
    23 = @{ D = "9.22";      E = "  -2.96%  " }
    24 = @{ E = "  -0.31%  " }
    25 = @{ D = "153.08";    E = "  -0.10%  " }
    26 = @{ E = "  +0.33%  " }
    27 = @{ D = "14.94";     E = "  -0.10%  " }
    28 = @{ E = "  +0.11%  " }
    29 = @{ E = "  -1.69%  " }
    30 = @{ D = "0.0463";    E = "  -1.90%  " }
    31 = @{ E = "  -3.52%  " }
    32 = @{ D = "3.16";      E = "  -1.31%  " }
    33 = @{ D = "1.399.78";  E = "  +0.38%  " }
    34 = @{ D = "2.91";      E = "  -2.28%  " }
    35 = @{ D = "1.51";      E = "  -2.97%  " }
    36 = @{ E = "  -0.67%  " }
    37 = @{ D = "0.925";     E = "  -1.77%  " }
    38 = @{ D = "0.0163";    E = "  -2.72%  " }
    39 = @{ D = "0.523";     E = "  -1.73%  " }
    40 = @{ D = "0.811";     E = "  -1.69%  " }
    41 = @{ E = "  +0.10%  " }
    42 = @{ D = "0.990";     E = "  -0.43%  " }
    43 = @{ D = "1.77";      E = "  -2.67%  " }
    44 = @{ D = "5.30";      E = "  +0.69%  " }
    45 = @{ D = "2.18";      E = "  +0.03%  " }
    46 = @{ D = "62.84";     E = "  -2.49%  " }
    47 = @{ D = "1.697.41";  E = "  -0.70%  " }
    48 = @{ D = "85.94";     E = "  -0.04%  " }
    49 = @{ D = "0.0₇0984"; E = "  -1.32%  " }
    50 = @{ B = "Cronos";    C = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro";       D = "0.0500"; E = "  +1.15%  " }
    51 = @{ B = "Algorand";  C = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo";    D = "0.0944"; E = "  -1.15%  " }
}

foreach ($row in $changes.Keys) {
    $cols = $changes[$row]
    foreach ($col in $cols.Keys) {
        $cell = $ws.Range("$col$row")
        $value = $cols[$col]
        if ($col -eq "D") {
            # Prices in column D are plain text in the source data (e.g. thousands
            # separated by dots like "1.565.51", or small decimals). Prefix with an
            # apostrophe so Excel stores them as text instead of auto-converting
            # number-looking values (like "205.61") into real numbers, then restore
            # the default "Normal" style since the text-entry coercion otherwise
            # tags the cell with an explicit (Text) number format style.
            $cell.Value = "'" + $value
            $cell.Style = "Normal"
        } else {
            $cell.Value = $value
        }
    }
}
